$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.799.50"
$ws.Range("D3").Value = "1.781.61"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "310.82"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D7").Value = "0.5121"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "0.3778"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").Value = "0.07762"
$ws.Range("E9").Value = "  -8.16%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "1.086"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "6.204"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "20.14"
$ws.Range("E14").Value = "  -4.21%  "
$ws.Range("D15").Value = "1.774.96"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "7.166"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("D17").Value = "91.98"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "0.00001071"
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("D19").Value = "0.06540"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D21").Value = "16.96"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").Value = "5.920"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "27.835.55"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "10.96"
$ws.Range("E24").Value = "  -4.10%  "
$ws.Range("D25").Value = "2.239"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "158.52"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  -4.30%  "
$ws.Range("D28").Value = "1.983.69"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "2.351"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "124.62"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "0.1077"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "1.026"
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("D33").Value = "3.615"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").Value = "5.486"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("D35").Value = "0.07055"
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("D36").Value = "0.02305"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").Value = "8.713"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "0.2118"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").Value = "11.50"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").Value = "5.010"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("D41").Value = "0.6085"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("E44").Value = "  -6.00%  "
$ws.Range("D45").Value = "0.5951"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "13.03"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").Value = "3.715"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "127.64"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "1.212"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "1.896"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("D51").Value = "0.06704"
$ws.Range("E51").Value = "  -3.98%  "
